# Updating the test data for PROD environment
# - Insert a new column I ("ExpectedSourceTemplateFile") holding the expected
#   source template path, shifting the old I/J/K ("ExpectedFilenames",
#   "StudyDesignExpectedValue", "ReportedVarExpectedValue") data right into
#   J/K/L.
# - Rename the "Pfizer - MM Maintenance" / "Interventional" population/slrtype
#   test row to "Takeda - MM Maintenance" / "Clinical".
# - Update the dependent Excel/Word report filename strings that referenced
#   "Pfizer - MM RRMM-...-Interventional-" / "...Maintenance-Interventional-"
#   so they read "Takeda - MM Maintenance-Clinical-" (and the related
#   Economic / Quality of Life / Real-world Evidence variants).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert the new column before the old column I; this shifts the existing
# I/J/K columns (and their column-width formatting) one slot to the right,
# matching the dimension growing from K31 to L31.
$ws.Columns("I:I").Insert()

# New column I: header + the single populated data value in row 2.
$ws.Range("I1").Value = "ExpectedSourceTemplateFile"
$ws.Range("I2").Value = "\Testdata\Templates\SLRReport_SourceData\PROD_Data\Expected_Source_Data_Manipulated.xlsx"

# Approximate the new column's (non bestFit) width from the diff as closely
# as this host's column-width quantization allows.
$ws.Columns("I:I").ColumnWidth = 27.5

# Rename the scenario row (A2:D2) from Pfizer/Interventional to
# Takeda/Clinical.
$ws.Range("A2").Value = "Takeda - MM Maintenance"
$ws.Range("B2").Value = "Takeda - MM Maintenance_radio_button"
$ws.Range("C2").Value = "Clinical"
$ws.Range("D2").Value = "Clinical_radio_button"

# Update the shifted ExpectedFilenames column (now column J) report-name
# strings to reference Takeda - MM Maintenance - Clinical (rows 3-10), and
# rename the remaining Pfizer Maintenance-Interventional rows to
# Pfizer Maintenance-Clinical (rows 11-12).
$ws.Range("J3").Value = "ExcelReport-Takeda - MM Maintenance-Clinical-"
$ws.Range("J4").Value = "WordReport-Takeda - MM Maintenance-Clinical-"
$ws.Range("J5").Value = "ExcelReport-Takeda - MM Maintenance-Economic-"
$ws.Range("J6").Value = "WordReport-Takeda - MM Maintenance-Economic-"
$ws.Range("J7").Value = "ExcelReport-Takeda - MM Maintenance-Quality of Life-"
$ws.Range("J8").Value = "WordReport-Takeda - MM Maintenance-Quality of Life-"
$ws.Range("J9").Value = "ExcelReport-Takeda - MM Maintenance-Real-world Evidence-"
$ws.Range("J10").Value = "WordReport-Takeda - MM Maintenance-Real-world Evidence-"
$ws.Range("J11").Value = "ExcelReport-Pfizer - MM Maintenance-Clinical-"
$ws.Range("J12").Value = "WordReport-Pfizer - MM Maintenance-Clinical-"

# Restore the active selection to I2 (matches the diff's updated
# <selection activeCell="I2" sqref="I2"/>).
$ws.Range("I2").Select()
